# Clear the "ragged edge" forecast-horizon cells that have no data.
# These were filtered out to support the new Component Analysis filtering options.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @(
    "G2", "H2", "I2", "J2", "K2",
    "I3", "J3", "K3",
    "G4", "H4", "I4", "J4", "K4",
    "I5", "J5", "K5",
    "G6", "H6", "I6", "J6", "K6",
    "I7", "J7", "K7",
    "G8", "H8", "I8", "J8", "K8",
    "I9", "J9", "K9",
    "G10", "H10", "I10", "J10", "K10",
    "I11", "J11", "K11",
    "G12", "H12", "I12", "J12", "K12",
    "I13", "J13", "K13",
    "G14", "H14", "I14", "J14", "K14",
    "I15", "J15", "K15",
    "G16", "H16", "I16", "J16", "K16",
    "I17", "J17", "K17",
    "K18",
    "I19", "J19", "K19",
    "K20",
    "I21", "J21", "K21",
    "K22",
    "J23", "K23",
    "I24", "J24", "K24",
    "K26",
    "J27", "K27",
    "I28", "J28", "K28",
    "K30",
    "J31", "K31",
    "I32", "J32", "K32",
    "K34",
    "J35", "K35",
    "I36", "J36", "K36",
    "K38",
    "J39", "K39",
    "I40", "J40", "K40",
    "K42",
    "J43", "K43",
    "I44", "J44"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

Write-Output ("Cleared {0} cells" -f $cellsToClear.Count)
